$d = $word.ActiveDocument

# 1. Weekly Progress Report: add ", individually submitted" before the closing paren.
$d.Content.Find.Execute(
    "Weekly Progress Report (during Project Period) ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Weekly Progress Report (during Project Period, individually submitted) ",
    2) | Out-Null

# 2. Mid-Term bullet: replace the due-date text with the paper announce/due info.
#    (The search text starts right after the "Mid-Term" hyperlink run, so it
#    does not need to span the hyperlink boundary.)
$d.Content.Find.Execute(
    ": Sunday 10/10/2021 11:59 PM (Fall Break: 10/11-10/12/2021)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " paper: Announced: Lesson 2; Due: Lesson 8.",
    2) | Out-Null
